$d = $word.ActiveDocument

# --- Update the placeholder ID text and drop the trailing run ---
# Original paragraph 1 has two runs: "**ID__AFFARS_mp_5315_3_topic_31__ID**"
# followed by a run that is just a single space (" "). The edit renames the
# placeholder (upper-cases "mp" -> "MP" and swaps the topic id for
# "5315_3_6_1") and removes the trailing space run entirely, leaving a
# single run with the new text and no extraneous whitespace.
$d.Content.Find.Execute(
    "**ID__AFFARS_mp_5315_3_topic_31__ID** ", $true, $false, $false, $false, $false,
    $true, 1, $false, "**ID__AFFARS_MP_5315_3_6_1__ID**", 2) | Out-Null

# --- Paragraph formatting for the (now single-run) first paragraph ---
$p1 = $d.Paragraphs(1)

# w:ind w:left="120" -> w:ind w:left="225"  (120 twips = 6pt, 225 twips = 11.25pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Add a paragraph border (w:pBdr) with 5-twip spacing on all four sides.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
